# Apply the "change adapt wce to fwce" edit to the workbook.
#
# Summary of changes:
#  1. Rename the two visible worksheets to append " Results".
#  2. Rename a set of "... Adaptive WCE" labels to "... FWCE" in column A
#     of the two test-result sheets (rows 9-17 and the "Median for ..."
#     summary rows).
#  3. Update the active sheet / selection so that "opus_big Test Results "
#     becomes the active tab (with cell A16 selected) while
#     "opus_base Test Results" keeps cell D36 selected (not the active tab).

$wb = $excel.ActiveWorkbook

# --- 1. Rename worksheets -------------------------------------------------
$wsBase = $wb.Worksheets.Item("opus_base Test")
$wsBase.Name = "opus_base Test Results"

$wsBig = $wb.Worksheets.Item("opus_big Test ")
$wsBig.Name = "opus_big Test Results "

# --- 2. Rename "Adaptive WCE" labels to "FWCE" -----------------------------
$renames = @{
    "Fine-Banded Adaptive WCE"               = "Fine-Banded FWCE"
    "Fine-Banded Unsampled Adaptive WCE"     = "Fine-Banded Unsampled FWCE"
    "Simple Adaptive WCE"                    = "Simple FWCE"
    "Simple Unsampled Adaptive WCE"          = "Simple Unsampled FWCE"
    "All-or-Nothing Adaptive WCE"            = "All-or-Nothing FWCE"
    "All-or-Nothing Unsampled Adaptive WCE"  = "All-or-Nothing Unsampled FWCE"
    "Fine-Banded LSP Adaptive WCE"           = "Fine-Banded LSP FWCE"
    "Simple Adaptive LSP WCE"                = "Simple LSP FWCE"
    "All-or-Nothing LSP Adaptive WCE"        = "All-or-Nothing LSP FWCE"
    "Median for adaptive WCE"                = "Median for FWCE"
    "Median for All-or-Nothing Adaptive WCE" = "Median for All-or-Nothing FWCE"
    "Median for Simple Adaptive WCE"         = "Median for Simple FWCE"
    "Median for Fine-Banded Adaptive WCE"    = "Median for Fine-Banded FWCE"
    "Median for Unsampled Adaptive WCE"      = "Median for Unsampled FWCE"
    "Median for Intersected Adaptive WCE"    = "Median for Intersected FWCE"
    "Median for LSP Adaptive WCE"            = "Median for LSP FWCE"
}

$sheetsToFix = @($wsBase, $wsBig)
foreach ($ws in $sheetsToFix) {
    # Column A rows 9-17 contain the WCE method name labels.
    for ($r = 9; $r -le 17; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($renames.ContainsKey($current)) {
            $cell.Value2 = $renames[$current]
        }
    }
    # The "Median for ..." summary labels live in column A somewhere in
    # rows 26-34 (exact row numbers differ slightly between the two
    # sheets because of a hidden row), so scan that range too.
    for ($r = 26; $r -le 34; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($renames.ContainsKey($current)) {
            $cell.Value2 = $renames[$current]
        }
    }
}

# --- 3. Update active sheet / selections -----------------------------------
$wsBase.Range("D36").Select()

$wsBig.Activate()
$wsBig.Range("A16").Select()
